$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.378.17"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'1.842.73"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'0.6274"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'24.86"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D12").Value = "'1.840.69"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'0.6740"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'0.00001023"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "'81.86"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'6.272"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "'29.365.40"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'234.39"
$ws.Range("D20").Value = "'12.33"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'7.302"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'157.68"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").Value = "'8.496"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'0.1345"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").Value = "'17.32"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "'0.07225"
$ws.Range("E28").Value = "  +10.85%  "
$ws.Range("D29").Value = "'1.489"
$ws.Range("E29").Value = "  +4.56%  "
$ws.Range("D30").Value = "'1.479"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.036"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.042"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "'1.819"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").Value = "'1.148"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "'0.6988"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "'2.573"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "'0.01834"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.797"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.840"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'1.233.52"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'0.9497"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'1.993.64"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "'100.98"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'65.25"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'1.704"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'6.964"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'8.886"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'0.1127"
$ws.Range("E51").Value = "  -2.13%  "
